# Update "Fresh bloom Flowers_2025-10-7.xlsx" workbook:
# - Append 20 new order-line rows (62-81) to the "Orders" sheet
# - Extend the concatenated "TotalNumber" digit-string in G2 of the "Summary" sheet

$wb = $excel.ActiveWorkbook
$ordersWs = $wb.Worksheets.Item("Orders")
$summaryWs = $wb.Worksheets.Item("Summary")

# Each tuple: (row, PackageID(A), FlowerName(C), Number(F)); $null = leave cell blank
$newRows = @(
    @(62, "15", "234_白泡泡_White Bubbles_Rosa rugosa Thunb._10stems", "9"),
    @(63, $null, "274_仙子之吻_undefined_Rosa rugosa Thunb._10stems", "5"),
    @(64, $null, "624_多丁白_undefined_undefined_1bunch", "5"),
    @(65, $null, "510_翠珠白_Didiscus caeruleus `nwhite_Trachymene Coerulea_1bunch", "10"),
    @(66, $null, "424_鼠尾白色_veronica white_undefined_1bunch", "5"),
    @(67, $null, "647_海棠果红_undefined_undefined_1bunch", "5"),
    @(68, "16", "321_雪柳叶_Spiraea  leaves_undefined_1bunch", "25"),
    @(69, $null, "320_雪柳花_Spiraea flower white_undefined_1bunch", "5"),
    @(70, $null, "542_吊米 红_hanging amaranthus`nred_undefined_1bunch", "5"),
    @(71, $null, "448_吊米 绿_hanging amaranthus`ngreen_undefined_1bunch", "3"),
    @(72, $null, "322_喷泉草_Grasses Panicum_undefined_1bunch", "5"),
    @(73, $null, "401_大飞燕白色_delphinium white_undefined_1bunch", "25"),
    @(74, $null, "495_大飞燕深粉色_delphinium pink_undefined_1bunch", "15"),
    @(75, $null, "411_紫罗兰白_violet white_undefined_1bunch", "20"),
    @(76, $null, "412_紫罗兰粉_violet pink_undefined_1bunch", "5"),
    @(77, $null, "319_尤加利叶带果_Eucalyptus leaves with small pods_undefined_1bunch", "5"),
    @(78, $null, "600_康乃馨复古红_vintage red_undefined_20stems", "8"),
    @(79, $null, "300_白星_White Gypso_ gypsophila_1kg", "9"),
    @(80, $null, "586_洋牡丹白_undefined_undefined_1bunch", "15"),
    @(81, $null, "585_洋牡丹红_undefined_undefined_1bunch", $null)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $aVal = $row[1]
    $cVal = $row[2]
    $fVal = $row[3]

    if ($null -ne $aVal) {
        $cell = $ordersWs.Cells.Item($r, 1)
        $cell.NumberFormat = "@"
        $cell.Value = $aVal
    }
    if ($null -ne $cVal) {
        $cell = $ordersWs.Cells.Item($r, 3)
        $cell.Value = $cVal
    }
    if ($null -ne $fVal) {
        $cell = $ordersWs.Cells.Item($r, 6)
        $cell.NumberFormat = "@"
        $cell.Value = $fVal
    }
}

# Extend the Summary sheet's concatenated-numbers text cell (G2) to include the new rows' counts
$g2 = $summaryWs.Range("G2")
$g2.NumberFormat = "@"
$g2.Value = "013242011.511125355101382212143175124050332553155521109102932355555555451252010205153015952055595510552555352515205589150"
